# Applies the crypto price/volume update described in the commit diff.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h) -- all stored as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.488.58'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.616.62'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.14'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.524'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.92'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.846.78'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.621.95'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.02'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.31'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.497.46'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '228.74'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0719'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.28'
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.87'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('E24').Value = '  +6.91%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '149.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.84'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.54'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0481'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.446.19'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.06'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.561'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.924'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.63%  '
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.860'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '68.93'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +6.34%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.39'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.21'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.757.54'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.67'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '86.09'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0982'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.53%  '
